$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 100.85
$ws.Range("I15").Value = 100.85
$ws.Range("K15").Value = 302.55
$ws.Range("M15").Value = -133.55
$ws.Range("H112").Value = 1314.7368
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 1317.5
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 3952.5
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -6168.5
$ws.Range("H129").Value = 646.2632
$ws.Range("J129").Value = 1006.6667
$ws.Range("L129").Value = 3020.0001
$ws.Range("N129").Value = -13020.0001
$ws.Range("H137").Value = 1308.6562
$ws.Range("I137").Value = 1327.4762
$ws.Range("J137").Value = 1272.7273
$ws.Range("K137").Value = 3982.4286
$ws.Range("L137").Value = 3818.1819
$ws.Range("M137").Value = -1432.4286
$ws.Range("N137").Value = -8918.1819
$ws.Range("H138").Value = 7069.581
$ws.Range("I138").Value = 3599
$ws.Range("J138").Value = 7951.932
$ws.Range("K138").Value = 10797
$ws.Range("L138").Value = 23855.796
$ws.Range("M138").Value = -5657
$ws.Range("N138").Value = -34135.796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 918.93335
$ws.Range("I2").Value = 1056.6364
$ws.Range("J2").Value = 540.25
$ws.Range("K2").Value = 1056.6364
$ws.Range("L2").Value = 540.25
$ws.Range("M2").Value = -943.6364000000001
$ws.Range("N2").Value = -766.25
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 40002
$ws.Range("I6").Value = 40002
$ws.Range("K6").Value = 40002
$ws.Range("M6").Value = -39829
$ws.Range("H24").Value = 20355
$ws.Range("J24").Value = 20355
$ws.Range("L24").Value = 20355
$ws.Range("N24").Value = -21103
$ws.Range("H32").Value = 18146.73
$ws.Range("I32").Value = 14937.632
$ws.Range("J32").Value = 26857.143
$ws.Range("K32").Value = 14937.632
$ws.Range("L32").Value = 26857.143
$ws.Range("M32").Value = -14650.632
$ws.Range("N32").Value = -27431.143
$ws.Range("H61").Value = 3002.182
$ws.Range("I61").Value = 2922.4
$ws.Range("J61").Value = 3800
$ws.Range("K61").Value = 2922.4
$ws.Range("L61").Value = 3800
$ws.Range("M61").Value = -2710.4
$ws.Range("N61").Value = -4224
$ws.Range("H100").Value = 20355
$ws.Range("J100").Value = 20355
$ws.Range("L100").Value = 20355
$ws.Range("N100").Value = -22519
$ws.Range("H116").Value = 918.93335
$ws.Range("I116").Value = 1056.6364
$ws.Range("J116").Value = 540.25
$ws.Range("K116").Value = 1056.6364
$ws.Range("L116").Value = 540.25
$ws.Range("M116").Value = 1237.3636
$ws.Range("N116").Value = -5128.25
$ws.Range("H122").Value = 2419
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850
$ws.Range("H132").Value = 2156.8086
$ws.Range("I132").Value = 1706.1333
$ws.Range("J132").Value = 2952.1177
$ws.Range("K132").Value = 5118.3999
$ws.Range("L132").Value = 8856.3531
$ws.Range("M132").Value = -2588.3999
$ws.Range("N132").Value = -13916.3531
$ws.Range("H136").Value = 3002.182
$ws.Range("I136").Value = 2922.4
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 8767.200000000001
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -6217.200000000001
$ws.Range("N136").Value = -16500
$ws.Range("H140").Value = 55333.332
$ws.Range("J140").Value = 55333.332
$ws.Range("L140").Value = 55333.332
$ws.Range("N140").Value = -65693.33199999999
$ws.Range("H141").Value = 58381.332
$ws.Range("J141").Value = 58381.332
$ws.Range("L141").Value = 58381.332
$ws.Range("N141").Value = -68741.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 918.93335
$ws.Range("I3").Value = 1056.6364
$ws.Range("J3").Value = 540.25
$ws.Range("K3").Value = 1056.6364
$ws.Range("L3").Value = 540.25
$ws.Range("M3").Value = -942.6364000000001
$ws.Range("N3").Value = -768.25
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 1875
$ws.Range("K134").Value = 5625
$ws.Range("M134").Value = -3090
$ws.Range("H140").Value = 59863.332
$ws.Range("J140").Value = 59863.332
$ws.Range("L140").Value = 59863.332
$ws.Range("N140").Value = -70223.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4833.304
$ws.Range("I31").Value = 2058.1
$ws.Range("J31").Value = 23334.666
$ws.Range("K31").Value = 2058.1
$ws.Range("L31").Value = 23334.666
$ws.Range("M31").Value = -1763.1
$ws.Range("N31").Value = -23924.666
$ws.Range("H34").Value = 4833.304
$ws.Range("I34").Value = 2058.1
$ws.Range("J34").Value = 23334.666
$ws.Range("K34").Value = 2058.1
$ws.Range("L34").Value = 23334.666
$ws.Range("M34").Value = -1856.1
$ws.Range("N34").Value = -23738.666
$ws.Range("H134").Value = 3191.04
$ws.Range("I134").Value = 3207.652
$ws.Range("K134").Value = 9622.956
$ws.Range("M134").Value = -7087.956
$ws.Range("H138").Value = 40519.875
$ws.Range("J138").Value = 40519.875
$ws.Range("L138").Value = 40519.875
$ws.Range("N138").Value = -50799.875
$ws.Range("H139").Value = 37007.3
$ws.Range("J139").Value = 37007.3
$ws.Range("L139").Value = 37007.3
$ws.Range("N139").Value = -47287.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1405.9286
$ws.Range("I113").Value = 2333.1667
$ws.Range("J113").Value = 710.5
$ws.Range("K113").Value = 6999.500100000001
$ws.Range("L113").Value = 2131.5
$ws.Range("M113").Value = -4829.500100000001
$ws.Range("N113").Value = -6471.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1908.4878
$ws.Range("I126").Value = 1670.6
$ws.Range("J126").Value = 2135.0476
$ws.Range("K126").Value = 5011.799999999999
$ws.Range("L126").Value = 6405.1428
$ws.Range("M126").Value = -2541.799999999999
$ws.Range("N126").Value = -11345.1428
$ws.Range("H132").Value = 3453.889
$ws.Range("I132").Value = 3341.7144
$ws.Range("J132").Value = 3525.2727
$ws.Range("K132").Value = 10025.1432
$ws.Range("L132").Value = 10575.8181
$ws.Range("M132").Value = -7495.143199999999
$ws.Range("N132").Value = -15635.8181
$ws.Range("H140").Value = 39750
$ws.Range("J140").Value = 39750
$ws.Range("L140").Value = 39750
$ws.Range("N140").Value = -50110

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1050
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1066.6666
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1066.6666
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1656.6666
$ws.Range("H27").Value = 1050
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1066.6666
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1066.6666
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1280.6666
$ws.Range("H122").Value = 15877787
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -51800
$ws.Range("H127").Value = 57642.777
$ws.Range("J127").Value = 57642.777
$ws.Range("L127").Value = 57642.777
$ws.Range("N127").Value = -67562.777
$ws.Range("H133").Value = 40478
$ws.Range("J133").Value = 40478
$ws.Range("L133").Value = 40478
$ws.Range("N133").Value = -45538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 48649
$ws.Range("J64").Value = 48649
$ws.Range("L64").Value = 48649
$ws.Range("N64").Value = -49145
$ws.Range("H67").Value = 48649
$ws.Range("J67").Value = 48649
$ws.Range("L67").Value = 48649
$ws.Range("N67").Value = -50365
$ws.Range("H81").Value = 2766.7778
$ws.Range("I81").Value = 1001
$ws.Range("J81").Value = 2987.5
$ws.Range("K81").Value = 2002
$ws.Range("L81").Value = 5975
$ws.Range("M81").Value = -941
$ws.Range("N81").Value = -8097
$ws.Range("H84").Value = 2766.7778
$ws.Range("I84").Value = 1001
$ws.Range("J84").Value = 2987.5
$ws.Range("K84").Value = 10010
$ws.Range("L84").Value = 29875
$ws.Range("M84").Value = -4706
$ws.Range("N84").Value = -40483
$ws.Range("H113").Value = 292
$ws.Range("I113").Value = 262.5
$ws.Range("J113").Value = 374.6
$ws.Range("K113").Value = 787.5
$ws.Range("L113").Value = 1123.8
$ws.Range("M113").Value = 1382.5
$ws.Range("N113").Value = -5463.8
$ws.Range("H137").Value = 57409.89
$ws.Range("J137").Value = 57409.89
$ws.Range("L137").Value = 57409.89
$ws.Range("N137").Value = -67609.89
